# ProjectSchedule.xlsx update — technology committee meeting 7/12/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the worksheet tab
# ---------------------------------------------------------------------
$ws.Name = "Schedule"

# ---------------------------------------------------------------------
# 2) Normalize per-column formatting across the whole task table
#    (A:no-wrap center, B:no-wrap center, C:wrap center, D/E:date,
#     F:no-wrap center, G:wrap center, H:wrap top) before rewriting
#    the row contents, by copying the existing formats down/over the
#    full A3:H18 block.
# ---------------------------------------------------------------------
$ws.Range("A3:B3").Copy() | Out-Null
$ws.Range("A4:B18").PasteSpecial(-4122) | Out-Null

$ws.Range("C3").Copy() | Out-Null
$ws.Range("C3:C18").PasteSpecial(-4122) | Out-Null

$ws.Range("D3:E3").Copy() | Out-Null
$ws.Range("D3:E18").PasteSpecial(-4122) | Out-Null

$ws.Range("F3").Copy() | Out-Null
$ws.Range("F3:F18").PasteSpecial(-4122) | Out-Null

$ws.Range("G3").Copy() | Out-Null
$ws.Range("G3:G18").PasteSpecial(-4122) | Out-Null

$ws.Range("H3").Copy() | Out-Null
$ws.Range("H3:H18").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Rewrite the task table contents (rows 3-18)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "High"
$ws.Range("C3").Value = "Randy Sheinbein"
$ws.Range("D3").Value = 45118
$ws.Range("E3").Value = 45119
$ws.Range("F3").Value = "2 Hrs"
$ws.Range("G3").Value = "Complete"
$ws.Range("H3").Value = "Get Website - Purchase Hosting From Provider"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "High"
$ws.Range("C4").Value = "Task 1"
$ws.Range("D4").Value = 45485
$ws.Range("E4").Value = 45119
$ws.Range("F4").Value = "24 Hrs"
$ws.Range("G4").Value = "Complete"
$ws.Range("H4").Value = "Point Domain Name Server to URL`nInform Mark Sheppard of IP address`nThis will take Mark about 5 minutes to configure the DNS record, it may take up to 2 days to reach all domain name servers"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Medium "
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = 45119
$ws.Range("E5").Value = 45120
$ws.Range("F5").Value = "4 Hrs"
$ws.Range("G5").Value = "Complete"
$ws.Range("H5").Value = "Put holding web page in place`nWelcome to the BHHS Class of 1974 50th Reunion Website"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "High"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = 45119
$ws.Range("E6").Clear() | Out-Null
$ws.Range("F6").Value = "2 Hrs"
$ws.Range("G6").Value = "In Progress"
$ws.Range("H6").Value = "Set up privileged user logins"

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "High"
$ws.Range("C7").Value = ""
$ws.Range("D7").Clear() | Out-Null
$ws.Range("E7").Clear() | Out-Null
$ws.Range("F7").Value = "2 Hrs"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = "Buy Divi Theme from Elegant Themes"

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "High"
$ws.Range("C8").Value = "Task 5"
$ws.Range("D8").Clear() | Out-Null
$ws.Range("E8").Clear() | Out-Null
$ws.Range("F8").Value = "6 Hrs"
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = "Set up basic theme (black, white orange)"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "High"
$ws.Range("C9").Value = "Ruth Amir"
$ws.Range("D9").Clear() | Out-Null
$ws.Range("E9").Clear() | Out-Null
$ws.Range("F9").Value = "1 Week"
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = "Graphic Design for website"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "High"
$ws.Range("C10").Value = "Task 7"
$ws.Range("D10").Clear() | Out-Null
$ws.Range("E10").Clear() | Out-Null
$ws.Range("F10").Value = "8 Hrs"
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = "Apply Graphic Design to website "

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Medium "
$ws.Range("C11").Value = "Task 8"
$ws.Range("D11").Clear() | Out-Null
$ws.Range("E11").Clear() | Out-Null
$ws.Range("F11").Value = "16 Hrs"
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = "Set up committee pages"

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Medium "
$ws.Range("C12").Value = "Task 9"
$ws.Range("D12").Clear() | Out-Null
$ws.Range("E12").Clear() | Out-Null
$ws.Range("F12").Value = "8 Hrs"
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = "Create forms for committee meeting minutes"

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "High"
$ws.Range("C13").Value = "Events Committee"
$ws.Range("D13").Value = 45124
$ws.Range("E13").Value = 45170
$ws.Range("F13").Value = "8 weeks"
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = "Select Venue`nSelect Date"

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "High"
$ws.Range("C14").Value = "Task 11"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("G14").Value = ""
$ws.Range("H14").Value = "Final Graphic Design with Ad Copy"

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "High"
$ws.Range("C15").Value = "Task 11"
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = "Elementary School Reunion web pages"

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "High"
$ws.Range("C16").Value = "Task 12"
$ws.Range("D16").Value = 45292
$ws.Range("E16").Value = 44985
$ws.Range("F16").Value = "TBD"
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = "Set up ecommerce portion of website for purchasing tickets"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "High"
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = 45138
$ws.Range("E17").Clear() | Out-Null
$ws.Range("F17").Value = "8 Hrs"
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = "Create Process for Contact Spreadsheet Update"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "High"
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = 45230
$ws.Range("E18").Clear() | Out-Null
$ws.Range("F18").Value = "2 Hrs"
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = "Mail blast delivery system"

# F14/F15 must not exist at all (no Estimated Effort logged for those tasks)
$ws.Range("F14").Clear() | Out-Null
$ws.Range("F15").Clear() | Out-Null

# ---------------------------------------------------------------------
# 4) Apply wrap text to column C (Dependencies) for the task rows,
#    matching the new style used throughout the table.
# ---------------------------------------------------------------------
$ws.Range("C3:C18").WrapText = $true

# ---------------------------------------------------------------------
# 5) Explicit row heights that changed with the new content
# ---------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30

# ---------------------------------------------------------------------
# 6) Extend the G/H placeholder rows down through row 38
# ---------------------------------------------------------------------
$ws.Range("G32:H32").Copy() | Out-Null
$ws.Range("G33:H38").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 7) Selection / view state
# ---------------------------------------------------------------------
$ws.Range("A18").Select() | Out-Null
